# JointFilerSSNmoreThan9Error.xlsx - RAD Phase 3 update for Estate Tax test data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the two "Date" timestamps used by the RAD rows (B2 / B4).
$ws.Range("B2").Value = "Thu Jan 25 17:48:24 EST 2024"
$ws.Range("B4").Value = "Thu Jan 25 17:48:37 EST 2024"

# Append a new row for the "Estate Tax" test case.
$ws.Range("D5").Value = "New Tax Return Amount Due"
$ws.Range("E5").Value = "Estate Tax"

# Give D5 the same border/wrap formatting used by the rest of column D.
$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# Leave the selection on the newly-added cell, matching the saved workbook state.
$ws.Range("E5").Select()
